# Get origins from Localiser.
# Refresh the SnippetID values (column H) on the "Voice Lines - main" sheet
# to the newly generated values coming from the Localiser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$snippetIds = @{
    2  = "B7t5"
    3  = "B7t5"
    4  = "jd0f"
    5  = "1G55"
    6  = "386V"
    7  = "Rfvk"
    8  = "rVib"
    9  = "xgUi"
    10 = "11u1"
    11 = "G6G4"
    12 = "18Wn"
    13 = "6IAN"
    14 = "ejHV"
    15 = "8qLB"
    16 = "twHV"
    17 = "Sl7I"
    18 = "Sl7I"
    19 = "Sl7I"
    20 = "Sl7I"
    21 = "Sl7I"
    22 = "VnUX"
    23 = "BUhZ"
    24 = "BuRx"
    25 = "u7lc"
    26 = "k9fF"
    27 = "k9fF"
    28 = "0twZ"
    29 = "ONXm"
}

foreach ($row in $snippetIds.Keys) {
    $ws.Cells.Item($row, 8).Value = $snippetIds[$row]
}
